$wb = $excel.ActiveWorkbook

# --- Sheet1 updates ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B3").Value = "No"
$ws1.Range("A2").Value = "Back Up WAN Circuit Down"
$ws1.Range("B2").Value = "Yes"
$ws1.Range("A3").Value = "Access to Network"

$ws1.Columns("A:A").AutoFit() | Out-Null

$ws1.Range("C6").Select() | Out-Null

# --- Sheet2 creation ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "CatalogSearch"
$ws2.Range("B1").Value = "execute"
$ws2.Range("A2").Value = "Broken keyboard or Mouse"
$ws2.Range("B2").Value = "Yes"
$ws2.Range("A3").Value = "Email Password Reset"
$ws2.Range("B3").Value = "No"

$ws2.Columns("A:A").AutoFit() | Out-Null
$ws2.Columns("B:B").AutoFit() | Out-Null

$ws2.Range("I11").Select() | Out-Null

$ws1.Activate() | Out-Null
